# "Generate Report for Handback" -- mark the two localized files (zh-cn, de-de)
# as handed back / in sync with en-US, fill in the per-language "Latest Target
# File" / "Latest Handback File" / "Latest Handback DateTime" columns, and
# widen the columns that now hold longer text.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# Known external hyperlink targets (same repo commit already referenced by the
# existing "Source File Name" hyperlinks on column A).
$hrefMd1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/66600597a83b1112b9500398b9017861a869914f/e2e/5183f9b5-3ddc-489e-ba8b-54204736914a.md"
$hrefMd2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/66600597a83b1112b9500398b9017861a869914f/e2e/530290e6-5491-4005-a69b-99d51fd2293c.md"

$md1 = "5183f9b5-3ddc-489e-ba8b-54204736914a.md"
$md2 = "530290e6-5491-4005-a69b-99d51fd2293c.md"

# ---------------------------------------------------------------------------
# Overview sheet: refresh the per-language status cells and widen the two
# language summary columns (E = zh-cn, F = de-de).
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus
$overview.Columns.Item(5).ColumnWidth = 29.083333333333332
$overview.Columns.Item(6).ColumnWidth = 29.083333333333332

# ---------------------------------------------------------------------------
# zh-cn sheet: handback just completed.
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$zhcn.Range("I2").Value = $md1
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $hrefMd1, "", "", $md1)
$zhcn.Range("J2").Value = "5183f9b5-3ddc-489e-ba8b-54204736914a.dbc3f4ce707a48212474bcb427163bf823842acf.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-09-06 06:56:00"

$zhcn.Range("I3").Value = $md2
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $hrefMd2, "", "", $md2)
$zhcn.Range("J3").Value = "530290e6-5491-4005-a69b-99d51fd2293c.971c62a76b38411303d77818c36246996a026c25.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-09-06 06:56:00"

$zhcn.Columns.Item(3).ColumnWidth = 29.083333333333332
$zhcn.Columns.Item(9).ColumnWidth = 39.083333333333336
$zhcn.Columns.Item(10).ColumnWidth = 39.083333333333336

# ---------------------------------------------------------------------------
# de-de sheet: handback just completed.
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

$dede.Range("I2").Value = $md1
$dede.Hyperlinks.Add($dede.Range("I2"), $hrefMd1, "", "", $md1)
$dede.Range("J2").Value = "5183f9b5-3ddc-489e-ba8b-54204736914a.dbc3f4ce707a48212474bcb427163bf823842acf.de-de.xlf"
$dede.Range("K2").Value = "2016-09-06 06:56:18"

$dede.Range("I3").Value = $md2
$dede.Hyperlinks.Add($dede.Range("I3"), $hrefMd2, "", "", $md2)
$dede.Range("J3").Value = "530290e6-5491-4005-a69b-99d51fd2293c.971c62a76b38411303d77818c36246996a026c25.de-de.xlf"
$dede.Range("K3").Value = "2016-09-06 06:56:18"

$dede.Columns.Item(3).ColumnWidth = 29.083333333333332
$dede.Columns.Item(9).ColumnWidth = 39.083333333333336
$dede.Columns.Item(10).ColumnWidth = 39.083333333333336
